$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for columns B,D,E,F,G,H,I,K,L,M across rows 2..25
# (column C, J, N, O remain 0 / unchanged)
$newData = @(
    @(0.5332579521459593,0.04165017798536041,0.2296257539761362,1.125842357355864,0.9914331776683838,0.9955867276757004,1.302672202578215,0.4569920692926814,0.1756469897096622,0.1559747625650907),
    @(0.5161907166284152,0.04111306026143069,0.2297578485753311,1.116667221472056,0.9833350096654812,0.9967228563899795,1.311400116601803,0.3981808411910208,0.1689583613377295,0.1506173307572141),
    @(0.5059383394914079,0.04077871432305002,0.2298533937188978,1.11171347993163,0.9790028559401804,0.9978943221015868,1.317277965508413,0.361983814306825,0.1649460447981994,0.1474013908382545),
    @(0.501817873994014,0.04064133067080888,0.2298959672477805,1.10986556612368,0.9773981607796003,0.9984908486363935,1.31980372438289,0.3472116977558528,0.1633348087207338,0.1461094508330802),
    @(0.5011371557896922,0.04061844996970976,0.2299032565052244,1.109569030878276,0.9771413996517566,0.9985970970243301,1.320231006230198,0.3447575064163004,0.1630687039515379,0.1458960500696875),
    @(0.5058825362782216,0.04077686610126108,0.2298539531389663,1.111687867085891,0.9789805641949414,0.9979018846792655,1.317311500422317,0.361784679209876,0.1649242186160933,0.1473838919297954),
    @(0.5273262858086696,0.04146593069248894,0.2296683078289075,1.122537562331203,0.9885078921376191,0.9958800997181356,1.305573902364785,0.436732092370363,0.173321156022638,0.1541123065798722),
    @(0.5711640402330715,0.04278067660737861,0.2294185232952737,1.149218372888839,1.012285235987093,0.995677153213137,1.286673187137865,0.5830114418057804,0.1905366339437933,0.1678872259830904),
    @(0.6044455282485899,0.04372394983199257,0.2293043074169328,1.172133879037133,1.032884138009138,0.9978252720655831,1.275296283989377,0.690065655173612,0.2036418574777343,0.1783587406808209),
    @(0.6198163558669876,0.04414806731963949,0.229267326435405,1.183282329637876,1.042940077454205,0.9993022919011594,1.270665540034308,0.7386785353841958,0.2097031746808398,0.1831982241992591),
    @(0.6256697629877408,0.04430794489668344,0.2292554705093339,1.18760833765019,1.046846924476455,0.9999335291648492,1.268990322946834,0.7570744150374082,0.2120127522169071,0.1850416609389285),
    @(0.6244076741687081,0.04427354487003043,0.2292579284507896,1.186672010885033,1.046001112188378,0.9997943814247776,1.26934762613427,0.7531131057060634,0.2115147086325919,0.1846441636529832),
    @(0.6202972639311213,0.0441612351615106,0.2292663080239192,1.183636140750679,1.043259512545887,0.9993527824777573,1.270526148632129,0.7401922330130901,0.2098928989215381,0.1833496683301234),
    @(0.6177837798488497,0.04409234731323863,0.2292717203210541,1.181790175958895,1.041593090370782,0.9990916581298279,1.271258231219527,0.7322761542189653,0.2089013530252686,0.1825581601884849),
    @(0.6034456211056067,0.04369613163526509,0.2293070250885316,1.171419889804341,1.032240778164294,0.9977388072044704,1.275609876022948,0.6868869197794538,0.2032477380981419,0.1780439893864951),
    @(0.5947084843840855,0.0434517828500276,0.2293325151571253,1.165243657688251,1.026679213088386,0.9970369286634337,1.278418987474254,0.6590196896713962,0.1998049271981301,0.1752940770363303),
    @(0.5897048725188654,0.04331077152736995,0.2293485863428328,1.16175939281375,1.023544841420502,0.9966802701821251,1.280085979280614,0.6429830611968157,0.1978340990381184,0.1737195533089633),
    @(0.5880144846034625,0.04326294738463687,0.2293542701087663,1.16059137735347,1.022494664028613,0.9965675900951965,1.280659198221421,0.6375519466103867,0.1971684238901474,0.1731876779509243),
    @(0.5956363185028977,0.04347784274224509,0.2293296558112257,1.16589407432248,1.027264574016584,0.9971067756648466,1.278114646841324,0.6619870480499799,0.2001704492574277,0.1755860702365588),
    @(0.6215037037983961,0.04419424304159492,0.2292637884884741,1.184525016187152,1.044062100653861,0.9994805384756944,1.270177861559759,0.7439877553473195,0.2103688763360054,0.1837295997996051),
    @(0.6386005846701437,0.04465821286481741,0.2292332579617034,1.197309636424578,1.055616758430773,1.001451167796375,1.265447379688965,0.797505424966829,0.2171173995245965,0.1891149471873277),
    @(0.6294583070639419,0.0444109745985557,0.229248409110097,1.190430521981312,1.049396963838888,1.000361029870902,1.267930334362305,0.7689489578889095,0.2135079825184221,0.1862349439271185),
    @(0.5952167839771789,0.0434660627205723,0.2293309441100275,1.165599813726914,1.026999735976929,0.9970750518638027,1.27825207727367,0.6606455532386519,0.2000051702870849,0.1754540400521591),
    @(0.5591150837915109,0.0424289550505712,0.2294739016773422,1.141420011718893,1.005304965817231,0.9953289264458078,1.29134563235317,0.5435126337835925,0.185799160110264,0.1640988897789555)
)

$cols = @("B","D","E","F","G","H","I","K","L","M")

for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $newData[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$rowNum").Value = $rowVals[$j]
    }
}
